$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Mene": update id_suivant (column B) values for the looping/ordering
# fix, then extend the concatenation formula in column G to also emit the
# id_precedent / id_suivant (A / B) values.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mene")
$ws.Activate()

$ws.Range("B2").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 2
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("B22").Value = 2
$ws.Range("B26").Value = 4
$ws.Range("B27").Value = 5
$ws.Range("B28").Value = 6
$ws.Range("B29").Value = 7
$ws.Range("B30").Value = 8
$ws.Range("B31").Value = 9
$ws.Range("B32").Value = 1
$ws.Range("B34").Value = 1
$ws.Range("B36").Value = 1
$ws.Range("B38").Value = 1
$ws.Range("B40").Value = 1

# G2 is a standalone formula; G3:G41 form a shared-formula group anchored on
# G3, so rewrite the whole range in one go to keep that group intact.
$ws.Range("G2").Formula = '="(''"&C2&"'',''"&D2&"'',"&A2&","&B2&"),"'
$ws.Range("G3:G41").Formula = '="(''"&C3&"'',''"&D3&"'',"&A3&","&B3&"),"'

[void]$ws.Range("A2:A5").Select()

# ---------------------------------------------------------------------------
# Sheet "Dialogue": the view simply scrolled back / re-selected a cell.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Dialogue")
$ws1.Activate()
[void]$ws1.Range("A3").Select()

# Leave "Mene" as the active sheet/tab, matching the saved workbook state.
$ws.Activate()
